$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.347.99"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").Value = "2.057.91"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("E6").Value = "  +2.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.93"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.79%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +3.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("D13").Value = "2.360.92"
$ws.Range("E13").Value = "  +4.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.78%  "

$ws.Range("E16").Value = "  +3.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.14%  "

$ws.Range("D18").Value = "2.058.30"
$ws.Range("E18").Value = "  +4.33%  "

$ws.Range("D19").Value = "37.527.62"
$ws.Range("E19").Value = "  +3.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "227.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.68%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  +3.62%  "

$ws.Range("E26").Value = "  +1.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("E28").Value = "  +13.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.28%  "

$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.73%  "

$ws.Range("E33").Value = "  +3.86%  "

$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.57%  "

$ws.Range("E36").Value = "  +6.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.09%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +31.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0990"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.48%  "

$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.08%  "

$ws.Range("D45").Value = "1.473.43"
$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.70%  "

$ws.Range("E49").Value = "  +3.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.76%  "

$ws.Range("E51").Value = "  +2.13%  "
